$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the corrected values for H10 and I10
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 5

# Highlight the corrected cells with a yellow fill
$ws.Range("H10:I10").Interior.Color = 65535

# Move the active selection to I14
$ws.Range("I14").Select()
